$wb = $excel.ActiveWorkbook

# 1) Rename the first sheet: ServiceProviderList -> ProviderList
$wsList = $wb.Worksheets.Item(1)
$wsList.Name = "ProviderList"

$wsIssue = $wb.Worksheets.Item(2)

# 2) Add new provider row 19 on ProviderList sheet (before the header rename so
#    new shared strings land ahead of "VendorCode" in the string table, matching
#    the canonical save order produced by Excel)
$wsList.Range("A19").Value = "NN012"
$wsList.Range("B19").Value = "test provider"
$wsList.Range("C19").Value = "Live"
$wsList.Range("D19").Value = "June"
$wsList.Range("E19").Value = "Martha"
$wsList.Range("F19").Value = "White Listing"
$wsList.Range("G19").Value = "TBD"
$wsList.Range("H19").Value = "IFM"
$wsList.Range("I19").Value = "chech provider date"

# 3) Rename header "ServiceProviderCode" -> "VendorCode" on both sheets (A1)
$wsList.Range("A1").Value = "VendorCode"
$wsIssue.Range("A1").Value = "VendorCode"

# 4) Update selections to match the saved view state
$wsIssue.Range("B18").Select()
$wsList.Range("D21").Select()
